# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new value for sheet "展览" (column F)
$sheet1Updates = @{
    2  = 1015
    5  = 7467
    9  = 780
    15 = 2993
    21 = 437
    22 = 21
    24 = 202
    25 = 207
    26 = 236
    28 = 86
    29 = 219
    37 = 69
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Row -> new value for sheet "全部类型" (column F)
$sheet4Updates = @{
    3  = 1015
    8  = 7467
    12 = 780
    19 = 2993
    27 = 437
    28 = 21
    30 = 202
    31 = 207
    32 = 236
    34 = 86
    35 = 219
    43 = 69
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
